# TimeTracking.xlsx update
#   "created DAO for tour, changed tour model"
#
# Adds the next batch of tour-log entries to the effort/time-tracking
# sheet (Tabelle1) and leaves the selection where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Extend the data-row formatting (date style in col A, hours style in
#     col C) from the first data row (row 4) down through row 40 so the
#     newly used rows pick up the same look as the existing table. ---
$ws.Range("A4").Copy()
$ws.Range("A5:A40").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C4").Copy()
$ws.Range("C5:C40").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- New tour-log entries ---
$ws.Range("A5").Value = 44676
$ws.Range("B5").Value = "Basic UI layout"
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = 44679
$ws.Range("B6").Value = "Finish basic UI layout and data binding"
$ws.Range("C6").Value = 4

$ws.Range("A7").Value = 44680
$ws.Range("B7").Value = "Define models, apply layered architecture"
$ws.Range("C7").Value = 2

$ws.Range("A8").Value = 44683
$ws.Range("B8").Value = "View Model Base, Relay Command"
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = 44685
$ws.Range("B9").Value = "Created tour view model, basic business logic"
$ws.Range("C9").Value = 2

$ws.Range("A10").Value = 44690
$ws.Range("B10").Value = "Created database with demo data, database connection"
$ws.Range("C10").Value = 3

$ws.Range("A11").Value = 44691

# --- Leave the selection the way the author left it ---
$ws.Range("B11").Select() | Out-Null
